$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "SCRIPT/P01P01A/us0106.ssb"
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = " We\'ve finally found it!"
$ws.Range("D5").Value = " Наконец-то мы его нашли!"
$ws.Range("E5").Value = " Îàëïîåø-óï íú åãï îàšìé!"

$ws.Range("A6").Value = "SCRIPT/P01P01A/us3102.ssb"
$ws.Range("B6").Value = 44
$ws.Range("C6").Value = " [CS:P]Zero Isle[CR]!"
$ws.Range("D6").Value = " Мы нашли [CS:P]Нуль-Остров[CR]!"
$ws.Range("E6").Value = " Íú îàšìé [CS:P]Îôìû-Ïòóñïâ[CR]!"

$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 43.2

$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C6").Select()
